$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.582.83"
$ws.Range("E2").Value = "  -2.62%  "
$ws.Range("D3").Value = "1.807.69"
$ws.Range("E4").Value = "  +0.49%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.67%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("E7").Value = "  +0.55%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "39.06"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -10.51%  "
$ws.Range("E10").Value = "  -3.39%  "
$ws.Range("E11").Value = "  -1.91%  "
$ws.Range("D12").Value = "2.069.01"
$ws.Range("E12").Value = "  -1.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.17"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.661"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.03%  "
$ws.Range("D15").Value = "1.791.41"
$ws.Range("E15").Value = "  -2.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.96%  "
$ws.Range("D17").Value = "34.606.68"
$ws.Range("E17").Value = "  -2.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("E20").Value = "  -2.50%  "
$ws.Range("E21").Value = "  -1.79%  "
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("E23").Value = "  +0.47%  "
$ws.Range("E24").Value = "  +1.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("E26").Value = "  -2.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.121"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("E29").Value = "  -7.17%  "
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("E31").Value = "  +2.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0543"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.90"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.50%  "
$ws.Range("E34").Value = "  +11.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.694"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "91.44"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.27%  "
$ws.Range("E38").Value = "  +4.62%  "
$ws.Range("D39").Value = "1.320.66"
$ws.Range("E39").Value = "  -2.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0191"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.03%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.955"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "14.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.27%  "
$ws.Range("E44").Value = "  -3.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.54%  "
$ws.Range("E46").Value = "  -0.85%  "
$ws.Range("E47").Value = "  -1.40%  "
$ws.Range("D48").Value = "1.992.47"
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("E49").Value = "  +0.49%  "
$ws.Range("E50").Value = "  +5.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "97.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.38%  "
